# Generate Report for Handoff
# Appends a new handed-off file ("98c3a135-a64a-4354-a872-aebd4fffea26")
# as row 3 on each of the three sheets (Overview, zh-cn, de-de), mirroring
# the existing row 2 entry for "46ac0aa9-8ff0-4dbb-aaf0-d0f3fd9d8efc".

$wb = $excel.ActiveWorkbook

$guid = "98c3a135-a64a-4354-a872-aebd4fffea26"
$hash = "0745cf5a6c178de4cfaa5f228efdf76f0693ab38"

$mdFile    = "$guid.md"
$zhCnXlf   = "$guid.$hash.zh-cn.xlf"
$deDeXlf   = "$guid.$hash.de-de.xlf"

$mdUrl   = "https://github.com/OpenLocalizationTest/oltest/blob/66cd11d8ea36d84bdd12b07dba03b56bbd7f6c98/e2e/$mdFile"
$zhCnUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/604b787e90d488c138ae56165cd5d720ddfb143e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$zhCnXlf"
$deDeUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/047de4d05593266c969521157be8c67f2e5d7c3a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$deDeXlf"

# ---------------------------------------------------------------------------
# Sheet "Overview": File Name | zh-cn | de-de | Latest Handoff Date
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A3").Value = $mdFile
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-27-12 10:27:14"

$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $mdUrl, "", "", $mdFile)

# ---------------------------------------------------------------------------
# Sheet "zh-cn": Source File Name | File Extension | Status | Latest Handoff
# File | Latest Handoff Datetime | Latest Target File | Latest Handback File
# | Latest Handback DateTime | Handoff Reason | Dependency From | Error Detail
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A3").Value = $mdFile
$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = $zhCnXlf
$wsZhCn.Range("E3").Value = "2016-03-12 10:27:11"
$wsZhCn.Range("H3").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("I3").Value = "Include"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $mdUrl, "", "", $mdFile)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("B3"), $mdUrl, "", "", ".md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D3"), $zhCnUrl, "", "", $zhCnXlf)

# ---------------------------------------------------------------------------
# Sheet "de-de": same columns as zh-cn
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A3").Value = $mdFile
$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = $deDeXlf
$wsDeDe.Range("E3").Value = "2016-03-12 10:27:14"
$wsDeDe.Range("H3").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("I3").Value = "Include"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $mdUrl, "", "", $mdFile)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("B3"), $mdUrl, "", "", ".md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D3"), $deDeUrl, "", "", $deDeXlf)
